$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 431
$ws.Range("I2").Value = 163.4
$ws.Range("J2").Value = 1100
$ws.Range("K2").Value = 163.4
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = -50.40000000000001
$ws.Range("N2").Value = -1326
$ws.Range("H6").Value = 75
$ws.Range("I6").Value = 75
$ws.Range("K6").Value = 225
$ws.Range("M6").Value = -113
$ws.Range("H19").Value = 641.1852
$ws.Range("I19").Value = 278.7143
$ws.Range("J19").Value = 768.05
$ws.Range("K19").Value = 278.7143
$ws.Range("L19").Value = 768.05
$ws.Range("M19").Value = -103.7143
$ws.Range("N19").Value = -1118.05
$ws.Range("H29").Value = 683.3333
$ws.Range("I29").Value = 683.3333
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2049.9999
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1768.9999
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 507.53845
$ws.Range("I38").Value = 276.44446
$ws.Range("J38").Value = 1027.5
$ws.Range("K38").Value = 829.33338
$ws.Range("L38").Value = 3082.5
$ws.Range("M38").Value = -457.33338
$ws.Range("N38").Value = -3826.5
$ws.Range("H53").Value = 245.28572
$ws.Range("I53").Value = 275.22223
$ws.Range("J53").Value = 222.83333
$ws.Range("K53").Value = 275.22223
$ws.Range("L53").Value = 222.83333
$ws.Range("M53").Value = 361.77777
$ws.Range("N53").Value = -1496.83333
$ws.Range("H58").Value = 1966.1111
$ws.Range("I58").Value = 739
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 2217
$ws.Range("L58").Value = 10500
$ws.Range("M58").Value = -2067
$ws.Range("N58").Value = -10800
$ws.Range("H87").Value = 45000
$ws.Range("J87").Value = 45000
$ws.Range("L87").Value = 45000
$ws.Range("N87").Value = -47496
$ws.Range("H90").Value = 45000
$ws.Range("J90").Value = 45000
$ws.Range("L90").Value = 135000
$ws.Range("N90").Value = -147480
$ws.Range("H98").Value = 1550.5
$ws.Range("I98").Value = 1474.8334
$ws.Range("J98").Value = 1701.8334
$ws.Range("K98").Value = 1474.8334
$ws.Range("L98").Value = 1701.8334
$ws.Range("M98").Value = 23.16660000000002
$ws.Range("N98").Value = -4697.8334
$ws.Range("H111").Value = 6503.7144
$ws.Range("I111").Value = 8105.2
$ws.Range("K111").Value = 24315.6
$ws.Range("M111").Value = -21248.6
$ws.Range("H122").Value = 1550.5
$ws.Range("I122").Value = 1474.8334
$ws.Range("J122").Value = 1701.8334
$ws.Range("K122").Value = 4424.5002
$ws.Range("L122").Value = 5105.5002
$ws.Range("M122").Value = -1974.5002
$ws.Range("N122").Value = -10005.5002
$ws.Range("H125").Value = 11664.2
$ws.Range("I125").Value = 2849.75
$ws.Range("J125").Value = 17540.5
$ws.Range("K125").Value = 25647.75
$ws.Range("L125").Value = 157864.5
$ws.Range("M125").Value = -23187.75
$ws.Range("N125").Value = -162784.5
$ws.Range("H135").Value = 46876532
$ws.Range("I135").Value = 20001394
$ws.Range("J135").Value = 142859170
$ws.Range("K135").Value = 180012546
$ws.Range("L135").Value = 1285732530
$ws.Range("M135").Value = -180010011
$ws.Range("N135").Value = -1285737600
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1781.2106
$ws.Range("I2").Value = 1865.6154
$ws.Range("J2").Value = 1598.3334
$ws.Range("K2").Value = 1865.6154
$ws.Range("L2").Value = 1598.3334
$ws.Range("M2").Value = -1752.6154
$ws.Range("N2").Value = -1824.3334
$ws.Range("H32").Value = 6013.5605
$ws.Range("I32").Value = 4788.629
$ws.Range("K32").Value = 4788.629
$ws.Range("M32").Value = -4501.629
$ws.Range("H61").Value = 5062.278
$ws.Range("I61").Value = 3436.4194
$ws.Range("J61").Value = 15142.6
$ws.Range("K61").Value = 3436.4194
$ws.Range("L61").Value = 15142.6
$ws.Range("M61").Value = -3224.4194
$ws.Range("N61").Value = -15566.6
$ws.Range("H63").Value = 1250.8334
$ws.Range("I63").Value = 1250.8334
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1250.8334
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -564.8334
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1250.8334
$ws.Range("I66").Value = 1250.8334
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 6254.166999999999
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -2822.166999999999
$ws.Range("N66").ClearContents()
$ws.Range("H116").Value = 1781.2106
$ws.Range("I116").Value = 1865.6154
$ws.Range("J116").Value = 1598.3334
$ws.Range("K116").Value = 1865.6154
$ws.Range("L116").Value = 1598.3334
$ws.Range("M116").Value = 428.3846000000001
$ws.Range("N116").Value = -6186.3334
$ws.Range("H136").Value = 5062.278
$ws.Range("I136").Value = 3436.4194
$ws.Range("J136").Value = 15142.6
$ws.Range("K136").Value = 10309.2582
$ws.Range("L136").Value = 45427.8
$ws.Range("M136").Value = -7759.2582
$ws.Range("N136").Value = -50527.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1781.2106
$ws.Range("I3").Value = 1865.6154
$ws.Range("J3").Value = 1598.3334
$ws.Range("K3").Value = 1865.6154
$ws.Range("L3").Value = 1598.3334
$ws.Range("M3").Value = -1751.6154
$ws.Range("N3").Value = -1826.3334
$ws.Range("H94").Value = 1546.9412
$ws.Range("I94").Value = 1518.0834
$ws.Range("J94").Value = 1616.2
$ws.Range("K94").Value = 1518.0834
$ws.Range("L94").Value = 1616.2
$ws.Range("M94").Value = -1067.0834
$ws.Range("N94").Value = -2518.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2140.7415
$ws.Range("I31").Value = 1426.721
$ws.Range("J31").Value = 4187.6
$ws.Range("K31").Value = 1426.721
$ws.Range("L31").Value = 4187.6
$ws.Range("N31").Value = -4777.6
$ws.Range("M31").Value = -1131.721
$ws.Range("H34").Value = 2140.7415
$ws.Range("I34").Value = 1426.721
$ws.Range("J34").Value = 4187.6
$ws.Range("K34").Value = 1426.721
$ws.Range("L34").Value = 4187.6
$ws.Range("N34").Value = -4591.6
$ws.Range("M34").Value = -1224.721
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 6782.7896
$ws.Range("I87").Value = 2711.75
$ws.Range("J87").Value = 9743.546
$ws.Range("K87").Value = 8135.25
$ws.Range("L87").Value = 29230.638
$ws.Range("M87").Value = -6887.25
$ws.Range("N87").Value = -31726.638
$ws.Range("H90").Value = 6782.7896
$ws.Range("I90").Value = 2711.75
$ws.Range("J90").Value = 9743.546
$ws.Range("K90").Value = 24405.75
$ws.Range("L90").Value = 87691.914
$ws.Range("M90").Value = -18165.75
$ws.Range("N90").Value = -100171.914
$ws.Range("H103").Value = 512.5
$ws.Range("I103").Value = 512.5
$ws.Range("K103").Value = 1537.5
$ws.Range("M103").Value = -658.5
$ws.Range("H114").Value = 620.8889
$ws.Range("I114").Value = 139.75
$ws.Range("J114").Value = 1005.8
$ws.Range("K114").Value = 419.25
$ws.Range("L114").Value = 3017.4
$ws.Range("M114").Value = 2834.75
$ws.Range("N114").Value = -9525.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 333361020
$ws.Range("I62").Value = 1000000000
$ws.Range("J62").Value = 41542.5
$ws.Range("K62").Value = 1000000000
$ws.Range("L62").Value = 41542.5
$ws.Range("N62").Value = -42914.5
$ws.Range("M62").Value = -999999314
$ws.Range("H63").Value = 36250
$ws.Range("J63").Value = 36250
$ws.Range("L63").Value = 36250
$ws.Range("N63").Value = -37622
$ws.Range("H65").Value = 333361020
$ws.Range("I65").Value = 1000000000
$ws.Range("J65").Value = 41542.5
$ws.Range("K65").Value = 3000000000
$ws.Range("L65").Value = 124627.5
$ws.Range("N65").Value = -131491.5
$ws.Range("M65").Value = -2999996568
$ws.Range("H66").Value = 36250
$ws.Range("J66").Value = 36250
$ws.Range("L66").Value = 108750
$ws.Range("N66").Value = -115614
$ws.Range("H107").Value = 408.83334
$ws.Range("I107").Value = 116.666664
$ws.Range("J107").Value = 701
$ws.Range("K107").Value = 116.666664
$ws.Range("L107").Value = 701
$ws.Range("M107").Value = 1803.333336
$ws.Range("N107").Value = -4541
$ws.Range("H122").Value = 5638.4614
$ws.Range("I122").Value = 8250
$ws.Range("J122").Value = 1460
$ws.Range("K122").Value = 24750
$ws.Range("L122").Value = 4380
$ws.Range("M122").Value = -22300
$ws.Range("N122").Value = -9280
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 758.1905
$ws.Range("I16").Value = 737.6111
$ws.Range("J16").Value = 881.6667
$ws.Range("K16").Value = 737.6111
$ws.Range("L16").Value = 881.6667
$ws.Range("M16").Value = -567.6111
$ws.Range("N16").Value = -1221.6667
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1282.3096
$ws.Range("I132").Value = 1356.8276
$ws.Range("J132").Value = 1116.0769
$ws.Range("K132").Value = 4070.4828
$ws.Range("L132").Value = 3348.2307
$ws.Range("M132").Value = -1540.4828
$ws.Range("N132").Value = -8408.2307
